$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $cellRef, $newValue)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" '28.788.72'
Set-TextCell $ws "E2" '  +3.23%  '
Set-TextCell $ws "D3" '1.880.65'
Set-TextCell $ws "E3" '  +3.23%  '
Set-TextCell $ws "E4" '  +0.33%  '
Set-TextCell $ws "D5" '324.94'
Set-TextCell $ws "E5" '  -1.03%  '
Set-TextCell $ws "D6" '1.005'
Set-TextCell $ws "E6" '  +0.34%  '
Set-TextCell $ws "D7" '0.4672'
Set-TextCell $ws "E7" '  +0.96%  '
Set-TextCell $ws "D8" '0.3937'
Set-TextCell $ws "E8" '  +2.47%  '
Set-TextCell $ws "D9" '0.07933'
Set-TextCell $ws "E9" '  +1.15%  '
Set-TextCell $ws "D10" '0.9777'
Set-TextCell $ws "E10" '  +2.07%  '
Set-TextCell $ws "D11" '22.35'
Set-TextCell $ws "E11" '  +2.35%  '
Set-TextCell $ws "D12" '1.892.76'
Set-TextCell $ws "E12" '  +1.32%  '
Set-TextCell $ws "D13" '5.748'
Set-TextCell $ws "E13" '  +1.90%  '
Set-TextCell $ws "D14" '7.015'
Set-TextCell $ws "E14" '  +2.43%  '
Set-TextCell $ws "D15" '0.06952'
Set-TextCell $ws "E15" '  +1.51%  '
Set-TextCell $ws "D16" '88.71'
Set-TextCell $ws "E16" '  +2.57%  '
Set-TextCell $ws "D17" '1.006'
Set-TextCell $ws "E17" '  +0.40%  '
Set-TextCell $ws "D18" '0.00001010'
Set-TextCell $ws "E18" '  +1.82%  '
Set-TextCell $ws "D19" '17.00'
Set-TextCell $ws "E19" '  +2.33%  '
Set-TextCell $ws "E20" '  +0.22%  '
Set-TextCell $ws "D21" '28.806.80'
Set-TextCell $ws "E21" '  +3.22%  '
Set-TextCell $ws "D22" '5.367'
Set-TextCell $ws "E22" '  +1.25%  '
Set-TextCell $ws "D23" '11.10'
Set-TextCell $ws "E23" '  +1.44%  '
Set-TextCell $ws "D24" '2.121'
Set-TextCell $ws "E24" '  +0.93%  '
Set-TextCell $ws "D25" '2.131.32'
Set-TextCell $ws "E25" '  +2.96%  '
Set-TextCell $ws "D26" '153.58'
Set-TextCell $ws "E26" '  +0.99%  '
Set-TextCell $ws "D27" '19.41'
Set-TextCell $ws "E27" '  +1.22%  '
Set-TextCell $ws "D28" '5.764'
Set-TextCell $ws "E28" '  +0.11%  '
Set-TextCell $ws "D29" '2.003'
Set-TextCell $ws "E29" '  +1.85%  '
Set-TextCell $ws "D30" '120.05'
Set-TextCell $ws "E30" '  +3.11%  '
Set-TextCell $ws "D31" '0.09381'
Set-TextCell $ws "E31" '  +1.85%  '
Set-TextCell $ws "D32" '0.9399'
Set-TextCell $ws "E32" '  +0.75%  '
Set-TextCell $ws "D33" '5.316'
Set-TextCell $ws "E33" '  +0.77%  '
Set-TextCell $ws "D34" '1.355'
Set-TextCell $ws "E34" '  +3.21%  '
Set-TextCell $ws "E35" '  +0.26%  '
Set-TextCell $ws "D36" '0.05930'
Set-TextCell $ws "E36" '  -0.07%  '
Set-TextCell $ws "D37" '0.02120'
Set-TextCell $ws "E37" '  -0.97%  '
Set-TextCell $ws "D38" '1.160'
Set-TextCell $ws "E38" '  +1.59%  '
Set-TextCell $ws "D39" '7.915'
Set-TextCell $ws "E39" '  +4.92%  '
Set-TextCell $ws "D40" '0.5721'
Set-TextCell $ws "E40" '  +2.82%  '
Set-TextCell $ws "D41" '0.1798'
Set-TextCell $ws "E41" '  +2.02%  '
Set-TextCell $ws "D42" '9.998'
Set-TextCell $ws "E42" '  +0.86%  '
Set-TextCell $ws "D43" '0.07295'
Set-TextCell $ws "E43" '  +4.40%  '
Set-TextCell $ws "D44" '11.90'
Set-TextCell $ws "E44" '  +3.11%  '
Set-TextCell $ws "D45" '0.5347'
Set-TextCell $ws "E45" '  +2.18%  '
Set-TextCell $ws "D46" '1.149'
Set-TextCell $ws "E46" '  -5.30%  '
Set-TextCell $ws "D47" '2.119'
Set-TextCell $ws "E47" '  -4.07%  '
Set-TextCell $ws "D48" '1.847'
Set-TextCell $ws "E48" '  +1.68%  '
Set-TextCell $ws "D49" '114.16'
Set-TextCell $ws "E49" '  +1.98%  '
Set-TextCell $ws "D50" '2.372'
Set-TextCell $ws "E50" '  +3.52%  '
Set-TextCell $ws "D51" '1.005'
Set-TextCell $ws "E51" '  +0.36%  '
